$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, old C->D)
$ws.Range("B1").EntireColumn.Insert()

# New column header (new model)
$ws.Range("B1").Value = "Sungrow SC2750UD-MV-US"

# Fill in new column B values for the new model
$ws.Range("B2").Value = 2750
$ws.Range("B3").Value = 2750
$ws.Range("B4").Value = 2750
$ws.Range("B5").Value = 2499.75
$ws.Range("B6").Value = 800
$ws.Range("B7").Value = 1500
$ws.Range("B8").Value = 250000
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 550
$ws.Range("B11").Formula = "=C11"
$ws.Range("B11").Locked = $True
$ws.Range("B11").NumberFormat = "0.00%"

# Selection as in the target file
$ws.Range("B12").Select()
